$d = $word.ActiveDocument

# 1. Patient name
$d.Content.Find.Execute("OTAVIO RAMOS DE ALMEIDA   ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "LUIZ CARLOS BOM   ", 2)

# 2. Birth date
$d.Content.Find.Execute("15/02/1988   ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "24/08/1953   ", 2)

# 3. Record number
$d.Content.Find.Execute("26294", $true, $false, $false, $false, $false,
                         $true, 1, $false, "121237", 2)

# 4. Mother's name
$d.Content.Find.Execute("MARIA IRACY RAMOS DOS SANTOS   ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "ROSALINA VOTORINO BOM   ", 2)

# 5. Date
$d.Content.Find.Execute("03/01/2019   ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "11/09/2018   ", 2)

# 6. Clear the evolution note text block, keeping an empty run in its place
#    (rather than letting the whole run be dropped, as happens when a run's
#    text is replaced down to "" in one step).
$needle = "# UROLOGIA" + [char]10 + "SOLICITO RETIRADA DE DUPLO J" + [char]10 + "30 PO DE URETERO + DUPLO J"
$full = $d.Content.Text
$start = $full.IndexOf($needle)
if ($start -ge 0) {
    $target = $d.Range($start, $start + $needle.Length)
    $target.Delete()
    $collapsed = $d.Range($start, $start)
    $collapsed.InsertAfter("")
}
